$wb = $excel.ActiveWorkbook

# --- 1. Fix the "Date Issued" number format on InvoiceLogTemplate!K6 ---
# It currently uses the custom "YYYY-MM-DD" (date-only) format; switch it
# back to the "YYYY-MM-DD HH:MM:SS" (date+time) format used elsewhere in
# the sheet.
$wsLog = $wb.Worksheets.Item("InvoiceLogTemplate")
$wsLog.Range("K6").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- 2. Append new client/project rows to the Clients sheet ---
$wsClients = $wb.Worksheets.Item("Clients")

$newRows = @(
    @("ATOS", "hello"),
    @("Cosmoco", "hello"),
    @("Cosmoco", "hello"),
    @("waqar", ""),
    @("waqar", "streamlit22"),
    @("ATOS", "wq"),
    @("ATOS", "wqqqq"),
    @("ATOS", "wqqqq12"),
    @("Minerva", "Machine Learning")
)

$startRow = 58
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $client = $newRows[$i][0]
    $project = $newRows[$i][1]
    $wsClients.Cells.Item($r, 2).Value = $client
    if ($project -ne "") {
        $wsClients.Cells.Item($r, 3).Value = $project
    }
}
